$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intro")

# New text added to the existing "parameters" row (row 51) describing usage.
$ws.Range("D51").Value = "Use with a row type of column, space separated parameters from the list below"

# Insert two new rows after the existing "rows=xx" row (row 52), shifting
# everything below down by two. Row 53 becomes a new "source=question"
# parameter row, row 54 is left blank (matching the spacing pattern used
# elsewhere in this reference sheet).
$ws.Rows("53:54").Insert()

$ws.Range("D53").Value = "Identify a question that will be the source of data for a column"
$ws.Range("C53").Value = "source=question"

$ws.Range("D53").Select()
